# Auto-generated script applying updated "F" column (想去人数 / interest count) values
# per the diff: update generated at commit 456a3b4
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 132  # was 133
$ws.Range("F4").Value = 1279  # was 1276
$ws.Range("F7").Value = 984  # was 981
$ws.Range("F8").Value = 943  # was 942
$ws.Range("F12").Value = 139  # was 38
$ws.Range("F14").Value = 932  # was 926
$ws.Range("F15").Value = 1812  # was 1807
$ws.Range("F16").Value = 4042  # was 3992
$ws.Range("F17").Value = 1183  # was 1169
$ws.Range("F18").Value = 113  # was 112
$ws.Range("F19").Value = 2637  # was 2620
$ws.Range("F21").Value = 1084  # was 1081
$ws.Range("F22").Value = 3617  # was 3598
$ws.Range("F23").Value = 769  # was 762
$ws.Range("F24").Value = 842  # was 840
$ws.Range("F25").Value = 42  # was 37
$ws.Range("F26").Value = 2329  # was 2321
$ws.Range("F27").Value = 114  # was 112
$ws.Range("F28").Value = 850  # was 846
$ws.Range("F29").Value = 169  # was 168
$ws.Range("F30").Value = 603  # was 556
$ws.Range("F33").Value = 1361  # was 1347
$ws.Range("F34").Value = 1965  # was 1957
$ws.Range("F36").Value = 497  # was 492
$ws.Range("F37").Value = 61  # was 51
$ws.Range("F39").Value = 592  # was 589
$ws.Range("F40").Value = 285  # was 283
$ws.Range("F41").Value = 77  # was 62
$ws.Range("F42").Value = 169  # was 167
$ws.Range("F43").Value = 236  # was 234
$ws.Range("F44").Value = 80  # was 78

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 2  # was 1
$ws.Range("F9").Value = 17  # was 16
$ws.Range("F12").Value = 119  # was 116
$ws.Range("F15").Value = 2  # was 1

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 445  # was 441

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 445  # was 441
$ws.Range("F3").Value = 132  # was 133
$ws.Range("F4").Value = 1279  # was 1276
$ws.Range("F6").Value = 984  # was 981
$ws.Range("F7").Value = 943  # was 942
$ws.Range("F14").Value = 932  # was 926
$ws.Range("F15").Value = 1812  # was 1807
$ws.Range("F16").Value = 4042  # was 3992
$ws.Range("F17").Value = 1183  # was 1169
$ws.Range("F18").Value = 113  # was 112
$ws.Range("F20").Value = 2637  # was 2620
$ws.Range("F21").Value = 1084  # was 1081
$ws.Range("F22").Value = 3617  # was 3598
$ws.Range("F23").Value = 769  # was 762
$ws.Range("F24").Value = 842  # was 840
$ws.Range("F26").Value = 42  # was 37
$ws.Range("F27").Value = 2329  # was 2321
$ws.Range("F28").Value = 17  # was 16
$ws.Range("F31").Value = 114  # was 112
$ws.Range("F32").Value = 119  # was 116
$ws.Range("F33").Value = 850  # was 846
$ws.Range("F34").Value = 169  # was 168
$ws.Range("F35").Value = 603  # was 557
$ws.Range("F38").Value = 1361  # was 1347
$ws.Range("F39").Value = 1965  # was 1957
$ws.Range("F40").Value = 2  # was 1
$ws.Range("F42").Value = 497  # was 492
$ws.Range("F43").Value = 61  # was 51
$ws.Range("F44").Value = 592  # was 589
$ws.Range("F45").Value = 285  # was 283
$ws.Range("F46").Value = 77  # was 62
$ws.Range("F47").Value = 169  # was 167
$ws.Range("F48").Value = 236  # was 234
$ws.Range("F49").Value = 80  # was 78
